$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 49  # was 48
$ws.Range("E4").Value = 22  # was 21
$ws.Range("E17").Value = 115  # was 113
$ws.Range("E18").Value = 107  # was 108
$ws.Range("E19").Value = 55  # was 54
$ws.Range("F19").Value = 28  # was 27
$ws.Range("H19").Value = 41  # was 40
$ws.Range("E33").Value = 41  # was 42
$ws.Range("G33").Value = 12  # was 13
$ws.Range("H33").Value = 23  # was 24
$ws.Range("E34").Value = 21  # was 20
$ws.Range("G34").Value = 3  # was 2
$ws.Range("H34").Value = 11  # was 10
$ws.Range("E36").Value = 96  # was 95
$ws.Range("F36").Value = 43  # was 42
$ws.Range("H36").Value = 75  # was 74
$ws.Range("E38").Value = 74  # was 73
$ws.Range("E41").Value = 39  # was 40
$ws.Range("F41").Value = 15  # was 16
$ws.Range("H41").Value = 26  # was 27
$ws.Range("E47").Value = 58  # was 57
$ws.Range("F47").Value = 36  # was 35
$ws.Range("H47").Value = 47  # was 46
$ws.Range("E49").Value = 70  # was 69
$ws.Range("E70").Value = 43  # was 42
$ws.Range("E72").Value = 41  # was 40
$ws.Range("E89").Value = 36  # was 34
